$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (e.g. "59.509.31"); force Text format so
# numeric-looking values are not auto-converted to numbers by COM.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.509.31"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "2.601.90"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "539.83"
$ws.Range("E5").Value = "  +3.34%  "

$ws.Range("D6").Value = "141.51"
$ws.Range("E6").Value = "  +1.35%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -0.49%  "

$ws.Range("E10").Value = "  +1.36%  "

$ws.Range("E11").Value = "  +1.40%  "

$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("D13").Value = "3.061.00"
$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").Value = "59.408.29"
$ws.Range("E14").Value = "  +0.78%  "

$ws.Range("D15").Value = "20.87"
$ws.Range("E15").Value = "  +1.37%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000133"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.571.47"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").Value = "341.30"
$ws.Range("E18").Value = "  +0.71%  "

$ws.Range("E19").Value = "  +1.02%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "6.36"
$ws.Range("E21").Value = "  -1.71%  "

$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "67.33"
$ws.Range("E23").Value = "  +1.76%  "

$ws.Range("E24").Value = "  +1.16%  "

$ws.Range("D25").Value = "0.166"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "7.24"
$ws.Range("E27").Value = "  +2.83%  "

$ws.Range("D28").Value = "0.0₃0744"
$ws.Range("E28").Value = "  +2.49%  "

$ws.Range("D30").Value = "1.67"
$ws.Range("E30").Value = "  +6.01%  "

$ws.Range("D31").Value = "5.83"
$ws.Range("E31").Value = "  -1.14%  "

$ws.Range("D32").Value = "18.81"
$ws.Range("E32").Value = "  +0.53%  "

$ws.Range("D33").Value = "149.83"
$ws.Range("E33").Value = "  +0.29%  "

$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("D35").Value = "1.12"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").Value = "0.837"
$ws.Range("E36").Value = "  +1.78%  "

$ws.Range("E37").Value = "  -0.78%  "

$ws.Range("D38").Value = "0.824"
$ws.Range("E38").Value = "  -0.55%  "

$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.20%  "

$ws.Range("D41").Value = "272.87"
$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").Value = "0.600"
$ws.Range("E42").Value = "  +1.44%  "

$ws.Range("D43").Value = "10.76"
$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("D44").Value = "0.0952"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("E45").Value = "  +1.40%  "

$ws.Range("D46").Value = "18.56"
$ws.Range("E46").Value = "  +3.49%  "

$ws.Range("E47").Value = "  +1.10%  "

$ws.Range("D48").Value = "1.941.81"
$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("D49").Value = "4.49"
$ws.Range("E49").Value = "  -0.75%  "

$ws.Range("D50").Value = "111.90"
$ws.Range("E50").Value = "  -1.18%  "

$ws.Range("E51").Value = "  +1.34%  "
